$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.733.16"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "'3.553.09"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'615.81"
$ws.Range("E5").Value = "  -5.85%  "
$ws.Range("D6").Value = "'154.31"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("D7").Value = "'3.549.97"
$ws.Range("E7").Value = "  -3.25%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.486"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "'6.90"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "'0.432"
$ws.Range("E12").Value = "  -1.74%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "'4.156.44"
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").Value = "'32.16"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'3.575.29"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "'67.749.57"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "'0.117"
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'15.59"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").Value = "'454.34"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").Value = "'9.47"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'77.71"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "'3.696.07"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -5.66%  "
$ws.Range("D28").Value = "'10.55"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  -6.43%  "
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "'25.95"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "'3.553.06"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("D38").Value = "'8.06"
$ws.Range("E38").Value = "  -3.90%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'177.18"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  -6.29%  "
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").Value = "'0.890"
$ws.Range("E45").Value = "  -4.27%  "
$ws.Range("D46").Value = "'28.99"
$ws.Range("E46").Value = "  +6.63%  "
$ws.Range("D47").Value = "'45.97"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("E48").Value = "  -4.76%  "
$ws.Range("D49").Value = "'7.68"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("E51").Value = "  -4.12%  "
